# rfq table broken links fix
#
# Updates the Purchase Request sheet: new PR number/date, replaces the
# broken "alternator and aircon belt" line item with the correct RFQ rows
# (Archfile Folder + Signpen), updates totals, and rewrites the Purpose
# text. Finally restores the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: PR No. and Date -------------------------------------------------
$ws.Range("C7").Value = "PR No.:  2020-03-0133"
$ws.Range("F7").Value = "March 02, 2020"

# --- Line item 1: Archfile Folder (replaces the old vehicle-repair line) ----
$ws.Range("A11").Value = "S273"
$ws.Range("B11").Value = "piece"
$ws.Range("C11").Value = "Archfile Folder, Legal" + [char]10 + "* 2`" /3`" spine 2 rings"
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 300
$ws.Range("F11").Value = 6000

# --- Line item 2: Signpen (new row, previously blank) -----------------------
$ws.Range("A12").Value = "S298"
$ws.Range("B12").Value = "piece"
$ws.Range("C12").Value = "Signpen, 0.7, Blue" + [char]10
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 100
$ws.Range("F12").Value = 500

# F36 holds =SUM(F11:F35) and recalculates automatically to 6500.

# --- Purpose -----------------------------------------------------------------
$ws.Range("B37").Value = "Realignment of Regional Records Management (QP-R4A-FAD-RICTU-08) to the newly developed Document Management System (DMS)"

# --- Restore the active cell selection ---------------------------------------
$ws.Range("E8").Select()
